# Updates cryptos list values (Price/Volume columns) and the
# Solana/XRP row swap (rows 6 and 7), per the Dec 20 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.891.26"
$ws.Range("E2").Value = "'  -0.50%  "
$ws.Range("D3").Value = "'2.212.52"
$ws.Range("E3").Value = "'  -1.22%  "
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("D5").Value = "'256.85"
$ws.Range("E5").Value = "'  +2.29%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = "'  +0.35%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'77.71"
$ws.Range("E7").Value = "'  +3.53%  "
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = "'  -0.79%  "
$ws.Range("D10").Value = "'43.02"
$ws.Range("E10").Value = "'  +5.12%  "
$ws.Range("D11").Value = "'0.0911"
$ws.Range("E11").Value = "'  -2.17%  "
$ws.Range("D12").Value = "'6.99"
$ws.Range("D13").Value = "'0.103"
$ws.Range("E13").Value = "'  +1.60%  "
$ws.Range("D14").Value = "'2.545.79"
$ws.Range("E14").Value = "'  -1.24%  "
$ws.Range("D15").Value = "'14.42"
$ws.Range("E15").Value = "'  -1.28%  "
$ws.Range("D16").Value = "'2.219.65"
$ws.Range("E16").Value = "'  -0.96%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "'  -1.20%  "
$ws.Range("D18").Value = "'42.840.44"
$ws.Range("E18").Value = "'  -0.37%  "
$ws.Range("E19").Value = "'  -0.92%  "
$ws.Range("D20").Value = "'71.04"
$ws.Range("E20").Value = "'  -0.15%  "
$ws.Range("E21").Value = "'  -0.17%  "
$ws.Range("D22").Value = "'2.30"
$ws.Range("E22").Value = "'  +4.69%  "
$ws.Range("D23").Value = "'229.82"
$ws.Range("E23").Value = "'  +0.04%  "
$ws.Range("D24").Value = "'9.25"
$ws.Range("E24").Value = "'  -4.94%  "
$ws.Range("E25").Value = "'  -0.14%  "
$ws.Range("D26").Value = "'42.88"
$ws.Range("E26").Value = "'  +8.74%  "
$ws.Range("D27").Value = "'10.72"
$ws.Range("E27").Value = "'  -0.53%  "
$ws.Range("D28").Value = "'3.35"
$ws.Range("E28").Value = "'  -2.64%  "
$ws.Range("E29").Value = "'  -2.72%  "
$ws.Range("D30").Value = "'2.20"
$ws.Range("E30").Value = "'  -0.94%  "
$ws.Range("D31").Value = "'173.44"
$ws.Range("E31").Value = "'  +1.10%  "
$ws.Range("D32").Value = "'20.37"
$ws.Range("E32").Value = "'  +0.73%  "
$ws.Range("D33").Value = "'0.0876"
$ws.Range("E33").Value = "'  +9.60%  "
$ws.Range("D34").Value = "'5.21"
$ws.Range("E34").Value = "'  -0.79%  "
$ws.Range("E35").Value = "'  -0.09%  "
$ws.Range("E36").Value = "'  +7.71%  "
$ws.Range("E37").Value = "'  -2.65%  "
$ws.Range("D38").Value = "'4.42"
$ws.Range("E38").Value = "'  -1.35%  "
$ws.Range("D39").Value = "'13.12"
$ws.Range("E39").Value = "'  +1.05%  "
$ws.Range("E40").Value = "'  +17.69%  "
$ws.Range("E41").Value = "'  -0.23%  "
$ws.Range("E42").Value = "'  -1.80%  "
$ws.Range("D43").Value = "'61.11"
$ws.Range("E43").Value = "'  +2.91%  "
$ws.Range("D44").Value = "'5.31"
$ws.Range("E44").Value = "'  -2.13%  "
$ws.Range("D45").Value = "'103.22"
$ws.Range("E45").Value = "'  -0.55%  "
$ws.Range("D46").Value = "'0.475"
$ws.Range("E46").Value = "'  -3.01%  "
$ws.Range("D47").Value = "'8.45"
$ws.Range("E47").Value = "'  -2.16%  "
$ws.Range("D48").Value = "'0.0971"
$ws.Range("E48").Value = "'  -1.92%  "
$ws.Range("E49").Value = "'  +1.25%  "
$ws.Range("E50").Value = "'  -1.45%  "
$ws.Range("E51").Value = "'  +22.57%  "
